$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.186.77'
$ws.Range('E2').Value = '  +2.20%  '
$ws.Range('D3').Value = '2.054.09'
$ws.Range('E3').Value = '  +1.64%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '231.87'
$ws.Range('E5').Value = '  +0.14%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.617'
$ws.Range('E6').Value = '  +3.22%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '57.35'
$ws.Range('E8').Value = '  +4.84%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.380'
$ws.Range('E9').Value = '  +3.05%  '
$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '57.87'
$ws.Range('E10').Value = '  +1.37%  '
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0755'
$ws.Range('E11').Value = '  +1.26%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.102'
$ws.Range('E12').Value = '  +1.29%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '2.360.16'
$ws.Range('E13').Value = '  +1.93%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '14.27'
$ws.Range('E14').Value = '  +0.43%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '20.70'
$ws.Range('E15').Value = '  +3.36%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.771'
$ws.Range('E16').Value = '  +1.77%  '
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '5.14'
$ws.Range('E17').Value = '  +0.93%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '2.060.52'
$ws.Range('E18').Value = '  +1.63%  '
$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = '37.135.12'
$ws.Range('E19').Value = '  +1.65%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.30'
$ws.Range('E20').Value = '  +14.63%  '
$ws.Range('B21').Value = 'Litecoin'
$ws.Range('C21').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '68.96'
$ws.Range('E21').Value = '  +2.27%  '
$ws.Range('B22').Value = 'ShibaInu'
$ws.Range('C22').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D22').Value = '0.0₃0807'
$ws.Range('E22').Value = '  +1.58%  '
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '224.10'
$ws.Range('E23').Value = '  +1.78%  '
$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.42'
$ws.Range('E25').Value = '  +1.76%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.36'
$ws.Range('E26').Value = '  +0.24%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '165.48'
$ws.Range('E27').Value = '  +1.91%  '
$ws.Range('B28').Value = 'ImmutableX'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.45'
$ws.Range('E28').Value = '  +7.20%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.74'
$ws.Range('E29').Value = '  +1.36%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '18.97'
$ws.Range('E30').Value = '  +0.65%  '
$ws.Range('B31').Value = 'Kaspa'
$ws.Range('C31').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.125'
$ws.Range('E31').Value = '  -2.42%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.117'
$ws.Range('E32').Value = '  +0.10%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.43'
$ws.Range('E33').Value = '  +1.67%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.53'
$ws.Range('E34').Value = '  +2.73%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0610'
$ws.Range('E35').Value = '  +1.55%  '
$ws.Range('B36').Value = 'InternetComputer(DFINITY)'
$ws.Range('C36').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.51'
$ws.Range('E36').Value = '  +6.20%  '
$ws.Range('B37').Value = 'BinanceUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('E38').Value = '  -0.96%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.26'
$ws.Range('E39').Value = '  -0.52%  '
$ws.Range('B40').Value = 'THORChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.73'
$ws.Range('E40').Value = '  -1.06%  '
$ws.Range('B41').Value = 'FTXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '4.63'
$ws.Range('E41').Value = '  +14.78%  '
$ws.Range('B42').Value = 'HuobiToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.95'
$ws.Range('E42').Value = '  +1.12%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '1.480.05'
$ws.Range('E43').Value = '  +0.63%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '95.79'
$ws.Range('E44').Value = '  +3.43%  '
$ws.Range('E45').Value = '  +4.74%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0926'
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0209'
$ws.Range('E47').Value = '  +3.33%  '
$ws.Range('E48').Value = '  +1.82%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '15.19'
$ws.Range('E49').Value = '  -2.54%  '
$ws.Range('B50').Value = 'FraxShare'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.11'
$ws.Range('E50').Value = '  +3.85%  '
$ws.Range('B51').Value = 'MXToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.94'
$ws.Range('E51').Value = '  +2.14%  '
